# prob25 (cal part) - add new key/buildingBlock rows 191-196, and shift the
# trailing z0001 marker row from 201 down to 211.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- row 191: existing A191 = "y0001" stays; add B191 / C191 ---
$ws.Range("B191").Value = '분모와 분자에 $n$을 곱해서 분모의 극한을 $0$ 아닌 값으로 확정시킵니다.'
$ws.Range("C191").Value = '32111_y23'

# --- row 192 ---
$ws.Range("A192").Value = 'y0002'
$ws.Range("B192").Value = '미분가능한 함수 $f$가 포함된 항등식의 양변을 $x$에 대해 미분해서 $f^{\prime}$이 포함된 항등식을 만듭니다.'
$ws.Range("C192").Value = '$f\left(x^{3}+x\right)=e^{x}$;'

# --- row 193 ---
$ws.Range("A193").Value = 'y0003'
$ws.Range("B193").Value = '양변에 적당한 $x$ 값을 대입해서 요구되는 미분계수를 구합니다.'

# --- row 194 ---
$ws.Range("A194").Value = 'y0004'
$ws.Range("B194").Value = '등비급수의 합을 구해서 첫째항과 공비의 방정식을 구합니다.'
$ws.Range("C194").Value = '$\displaystyle\sum_{n=1}^{\infty}\left(a_{2 n-1}-a_{2 n}\right)=3, \quad \displaystyle\sum_{n=1}^{\infty} a_{n}^{2}=6$'

# --- row 195 ---
$ws.Range("A195").Value = 'y0005'
$ws.Range("B195").Value = '두 등비급수에서 구한 첫째항과 공비의 연립방정식을 풉니다.'

# --- row 196 ---
$ws.Range("A196").Value = 'y0006'
$ws.Range("B196").Value = '등비급수의 합을 첫째항과 공비를 이용해서 구합니다.'

# --- move the trailing marker row from 201 to 211 ---
$ws.Range("A211").Value = $ws.Range("A201").Value2
$ws.Range("A201").ClearContents()

# --- restore view state (scroll position / selection) ---
$ws.Activate() | Out-Null
$excel.Goto($ws.Range("A184"), $true) | Out-Null
$ws.Range("B197").Select() | Out-Null
